$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new time-record row (row 3) for working on the VGA interface.
# Duplicate the date cell from row 2 (same day) via copy/paste so it keeps
# being stored as plain shared-string text instead of being reinterpreted
# as a date value/style.
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4163)  # xlPasteAll
$excel.CutCopyMode = 0

$ws.Range("B3").Value = 0.41666666666666669   # 10:00
$ws.Range("C3").Value = 0.46527777777777773   # 11:10
$ws.Range("D3").Formula = "=C3-B3"
$ws.Range("E3").Value = "VGA Interface"
$ws.Range("F3").Value = "Read into subject, take notes"

# Reflect the saved selection state (cell below the new row's Task column).
$ws.Range("E9").Select()
